$wb = $excel.ActiveWorkbook

# --- "Other settings" sheet: insert a new row for "Print html report" = "yes"
# (the existing row, previously "Print html report" = "no", is pushed down and
# re-labelled "Save GIS layers" since its text changes in the process) ---
$settings = $wb.Worksheets.Item("Other settings")
$settings.Rows.Item(7).Insert()
$settings.Range("A7").Value = "Print html report"
$settings.Range("B7").Value = "yes"
$settings.Range("A8").Value = "Save GIS layers"
$settings.Range("B8").Value = "no"

# extend the "yes/no" list validation down onto the newly inserted row
$settings.Range("B7:B8").Validation.Add(3, 1, 1, "tech!`$B`$2:`$B`$3")

# --- "Pest_list" sheet: replace the pest with a new one, and make it the
# active sheet/selection (matching the author's last interaction) ---
$pestList = $wb.Worksheets.Item("Pest_list")
$pestList.Range("A2").Value = "Amyelois transitella"
$pestList.Activate()
$pestList.Range("A3").Select()
